# Update the "想去人数" (interested count) figures in column F
# for the affected rows on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3020
    7  = 1653
    15 = 20
    21 = 3135
    22 = 386
    23 = 117
    24 = 188
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
